$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: product swap (SMART WATCH MI / XIAOMI -> SMART WATCH LG / LG) ---
$ws.Range("B19").Value = "SMART WATCH LG"
$ws.Range("C19").Value = "LG"

# --- New rows 26-28: Redmi phones ---
$ws.Cells.Item(26, 1).Value = "PD1025"
$ws.Cells.Item(26, 2).Value = "REDMI NOTE 7"
$ws.Cells.Item(26, 3).Value = "REDMI"
$ws.Cells.Item(26, 4).Value = "CELULAR"
$ws.Cells.Item(26, 5).Value = 1200
$ws.Cells.Item(26, 6).Value = 500

$ws.Cells.Item(27, 1).Value = "PD1026"
$ws.Cells.Item(27, 2).Value = "REDMI NOTE 8"
$ws.Cells.Item(27, 3).Value = "REDMI"
$ws.Cells.Item(27, 4).Value = "CELULAR"
$ws.Cells.Item(27, 5).Value = 1300
$ws.Cells.Item(27, 6).Value = 600

$ws.Cells.Item(28, 1).Value = "PD1027"
$ws.Cells.Item(28, 2).Value = "REDMI NOTE 9"
$ws.Cells.Item(28, 3).Value = "REDMI"
$ws.Cells.Item(28, 4).Value = "CELULAR"
$ws.Cells.Item(28, 5).Value = 1400
$ws.Cells.Item(28, 6).Value = 700

# --- AutoFilter range grows to A1:F28 (re-toggle since it was already on) ---
# Must happen here, while the used range still ends at row 28, otherwise Excel
# auto-expands the filter to cover rows added afterwards.
$ws.AutoFilterMode = $false
$ws.Range("A1:F28").AutoFilter()

# --- _FilterDatabase defined name follows the same range ---
$nm = $wb.Names.Item(1)
$nm.RefersTo = "=Planilha1!`$A`$1:`$F`$28"

# Row 29: Redmi Note 10 - E29 ends up with its own (bold-toggled/reverted) style, like the source file
$ws.Cells.Item(29, 1).Value = "PD1028"
$ws.Cells.Item(29, 2).Value = "REDMI NOTE 10"
$ws.Cells.Item(29, 3).Value = "REDMI"
$ws.Cells.Item(29, 4).Value = "CELULAR"
$ws.Cells.Item(29, 6).Value = 1400
$ws.Cells.Item(29, 5).Value = 2000
$ws.Range("E29").Font.Bold = $true
$ws.Range("E29").Font.Bold = $false
$ws.Range("E29").HorizontalAlignment = -4108

# --- New rows 30-31: Guitars ---
$ws.Cells.Item(30, 1).Value = "PD1029"
$ws.Cells.Item(30, 2).Value = "GUITARRA STRINBERG"
$ws.Cells.Item(30, 3).Value = "STRINBERG"
$ws.Cells.Item(30, 4).Value = "GUITARRA"
$ws.Cells.Item(30, 5).Value = 1800
$ws.Cells.Item(30, 6).Value = 1100

$ws.Cells.Item(31, 1).Value = "PD1030"
$ws.Cells.Item(31, 2).Value = "GUITARRA IBANEZ"
$ws.Cells.Item(31, 3).Value = "IBANEZ"
$ws.Cells.Item(31, 4).Value = "GUITARRA"
$ws.Cells.Item(31, 5).Value = 3000
$ws.Cells.Item(31, 6).Value = 2100

# --- Selection moved to F32, just past the new data ---
$ws.Range("F32").Select()

# --- Page setup: A4, portrait ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
